$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 44195
$ws.Range("K2").Value = 'Patterson'
$ws.Range("M2").Value = 124
$ws.Range("N2").Value = 13000
$ws.Range("O2").Value = 13000
$ws.Range("P2").Value = 13000
$ws.Range("Q2").Value = '$/caja 15 kilos'
$ws.Range("S2").Value = 867
$ws.Range("T2").Value = 15

$ws.Range("D3").Value = 44187
$ws.Range("K3").Value = 'Patterson'
$ws.Range("M3").Value = 80
$ws.Range("N3").Value = 15000
$ws.Range("P3").Value = 15000
$ws.Range("Q3").Value = '$/caja 15 kilos granel'
$ws.Range("S3").Value = 1000
$ws.Range("T3").Value = 15

$ws.Range("D4").Value = 44187
$ws.Range("L4").Value = 'Segunda'
$ws.Range("M4").Value = 95
$ws.Range("N4").Value = 13500
$ws.Range("O4").Value = 13500
$ws.Range("P4").Value = 13500
$ws.Range("Q4").Value = '$/caja 15 kilos granel'
$ws.Range("S4").Value = 900
$ws.Range("T4").Value = 15

$ws.Range("D5").Value = 44187
$ws.Range("K5").Value = 'Patterson'
$ws.Range("L5").Value = 'Tercera'
$ws.Range("M5").Value = 120
$ws.Range("N5").Value = 12000
$ws.Range("P5").Value = 12000
$ws.Range("Q5").Value = '$/caja 15 kilos granel'
$ws.Range("S5").Value = 800
$ws.Range("T5").Value = 15

$ws.Range("D6").Value = 44174
$ws.Range("K6").Value = 'Modesto'
$ws.Range("N6").Value = 8500
$ws.Range("O6").Value = 8500
$ws.Range("P6").Value = 8500
$ws.Range("Q6").Value = '$/bandeja 10 kilos'
$ws.Range("R6").Value = 'Región Metropolitana'
$ws.Range("S6").Value = 850
$ws.Range("T6").Value = 10

$ws.Range("D7").Value = 44174
$ws.Range("K7").Value = 'Modesto'
$ws.Range("M7").Value = 180
$ws.Range("N7").Value = 15000
$ws.Range("O7").Value = 15000
$ws.Range("P7").Value = 15000
$ws.Range("Q7").Value = '$/caja 18 kilos'
$ws.Range("R7").Value = 'Región Metropolitana'
$ws.Range("S7").Value = 833
$ws.Range("T7").Value = 18

$ws.Range("D8").Value = 44174
$ws.Range("K8").Value = 'Modesto'
$ws.Range("M8").Value = 120
$ws.Range("N8").Value = 12000
$ws.Range("O8").Value = 12000
$ws.Range("P8").Value = 12000
$ws.Range("Q8").Value = '$/caja 18 kilos'
$ws.Range("R8").Value = 'Región Metropolitana'
$ws.Range("S8").Value = 667
$ws.Range("T8").Value = 18

$ws.Range("D12").Value = 44176
$ws.Range("M12").Value = 115
$ws.Range("P12").Value = 11609
$ws.Range("S12").Value = 967

$ws.Range("D13").Value = 44162
$ws.Range("K13").Value = 'Castle Brite'
$ws.Range("M13").Value = 70
$ws.Range("R13").Value = 'Provincia de San Felipe de Aconcagua'

$ws.Range("D14").Value = 44162
$ws.Range("K14").Value = 'Castle Brite'
$ws.Range("M14").Value = 75
$ws.Range("N14").Value = 14000
$ws.Range("P14").Value = 14400
$ws.Range("R14").Value = 'Provincia de San Felipe de Aconcagua'
$ws.Range("S14").Value = 800

$ws.Range("D15").Value = 44194
$ws.Range("K15").Value = 'Patterson'
$ws.Range("L15").Value = 'Primera'
$ws.Range("N15").Value = 13000
$ws.Range("O15").Value = 13000
$ws.Range("P15").Value = 13000
$ws.Range("Q15").Value = '$/caja 15 kilos'
$ws.Range("R15").Value = 'Provincia de San Felipe de Aconcagua'
$ws.Range("S15").Value = 867
$ws.Range("T15").Value = 15

$ws.Range("D16").Value = 44159
$ws.Range("K16").Value = 'Castle Brite'
$ws.Range("M16").Value = 80
$ws.Range("N16").Value = 8000
$ws.Range("O16").Value = 8000
$ws.Range("P16").Value = 8000
$ws.Range("Q16").Value = '$/bandeja 10 kilos'
$ws.Range("S16").Value = 800
$ws.Range("T16").Value = 10

$ws.Range("D17").Value = 44159
$ws.Range("K17").Value = 'Castle Brite'
$ws.Range("L17").Value = 'Segunda'
$ws.Range("M17").Value = 65
$ws.Range("N17").Value = 7000
$ws.Range("O17").Value = 7000
$ws.Range("P17").Value = 7000
$ws.Range("R17").Value = 'Provincia de San Felipe de Aconcagua'
$ws.Range("S17").Value = 700

$ws.Range("D18").Value = 44169
$ws.Range("K18").Value = 'Dina'
$ws.Range("N18").Value = 10000
$ws.Range("O18").Value = 10000
$ws.Range("P18").Value = 10000
$ws.Range("R18").Value = 'Región de O''Higgins'
$ws.Range("S18").Value = 1000

$ws.Range("D19").Value = 44160
$ws.Range("L19").Value = 'Primera'
$ws.Range("M19").Value = 25
$ws.Range("N19").Value = 8000
$ws.Range("O19").Value = 8000
$ws.Range("P19").Value = 8000
$ws.Range("S19").Value = 800

$ws.Range("D20").Value = 44160
$ws.Range("L20").Value = 'Segunda'
$ws.Range("M20").Value = 40
$ws.Range("N20").Value = 7000
$ws.Range("O20").Value = 7000
$ws.Range("P20").Value = 7000
$ws.Range("R20").Value = 'Provincia de San Felipe de Aconcagua'
$ws.Range("S20").Value = 700

$ws.Range("D21").Value = 44189
$ws.Range("K21").Value = 'Patterson'
$ws.Range("L21").Value = 'Primera'
$ws.Range("M21").Value = 130
$ws.Range("N21").Value = 12000
$ws.Range("O21").Value = 12000
$ws.Range("P21").Value = 12000
$ws.Range("Q21").Value = '$/caja 18 kilos'
$ws.Range("R21").Value = 'Provincia de San Felipe de Aconcagua'
$ws.Range("S21").Value = 667
$ws.Range("T21").Value = 18

$ws.Range("D22").Value = 44166
$ws.Range("K22").Value = 'Castle Brite'
$ws.Range("M22").Value = 120
$ws.Range("N22").Value = 10000
$ws.Range("O22").Value = 10000
$ws.Range("P22").Value = 10000
$ws.Range("Q22").Value = '$/bandeja 10 kilos'
$ws.Range("R22").Value = 'Región Metropolitana'
$ws.Range("T22").Value = 10

$ws.Range("D23").Value = 44166
$ws.Range("K23").Value = 'Castle Brite'
$ws.Range("M23").Value = 120
$ws.Range("N23").Value = 8000
$ws.Range("O23").Value = 8000
$ws.Range("P23").Value = 8000
$ws.Range("Q23").Value = '$/bandeja 10 kilos'
$ws.Range("R23").Value = 'Región Metropolitana'
$ws.Range("S23").Value = 800
$ws.Range("T23").Value = 10

$ws.Range("D24").Value = 44175
$ws.Range("K24").Value = 'Modesto'
$ws.Range("L24").Value = 'Primera'
$ws.Range("M24").Value = 140
$ws.Range("N24").Value = 11000
$ws.Range("P24").Value = 11571
$ws.Range("Q24").Value = '$/caja 12 kilos'
$ws.Range("S24").Value = 967
$ws.Range("T24").Value = 12
